# Final Patch Update before release
#
# Adds English translations for the Yes/No choice strings that already
# live in column A:
#   Row 13 -> A13 = "はい"   (Yes)  => D13 = "Yes"
#   Row 14 -> A14 = "いいえ" (No)   => D14 = "No"
#
# This introduces a new column D (dimension grows from A1:C49 to
# A1:D49) without touching any of the existing data in columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "Yes"
$ws.Range("D14").Value = "No"
